# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F2").Value = 137
        $ws.Range("F5").Value = 6748
        $ws.Range("F7").Value = 4
        $ws.Range("F9").Value = 143
        $ws.Range("F10").Value = 6284
        $ws.Range("F11").Value = 51
        $ws.Range("F12").Value = 193
        $ws.Range("F13").Value = 1268
        $ws.Range("F17").Value = 124
        $ws.Range("F19").Value = 369
        $ws.Range("F20").Value = 46
        $ws.Range("F22").Value = 4631
        $ws.Range("F23").Value = 64
        $ws.Range("F24").Value = 43
        $ws.Range("F25").Value = 96
        $ws.Range("F27").Value = 77
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F2").Value = 137
        $ws.Range("F5").Value = 6748
        $ws.Range("F7").Value = 4
        $ws.Range("F9").Value = 143
        $ws.Range("F10").Value = 6284
        $ws.Range("F11").Value = 51
        $ws.Range("F12").Value = 193
        $ws.Range("F13").Value = 1268
        $ws.Range("F17").Value = 124
        $ws.Range("F19").Value = 369
        $ws.Range("F20").Value = 46
        $ws.Range("F22").Value = 4631
        $ws.Range("F24").Value = 64
        $ws.Range("F25").Value = 43
        $ws.Range("F26").Value = 96
        $ws.Range("F28").Value = 77
    }
}
